$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A92").Value = "11.01.2022 17:23 (Kyiv+Israel) 15:23 (UTC) 00:23 (Japan) 20:53 (India)"
$ws.Range("B92").Value = "***"
$ws.Range("C92").Value = "***"
$ws.Range("D92").Value = 6.476
$ws.Range("E92").Value = -5.338
$ws.Range("A93").Value = "11.01.2022 18:01 (Kyiv+Israel) 16:01 (UTC) 01:01 (Japan) 21:31 (India)"
$ws.Range("B93").Value = "***"
$ws.Range("C93").Value = "***"
$ws.Range("D93").Value = 3.201
$ws.Range("E93").Value = -2.063
$ws.Range("A94").Value = "11.02.2022 13:46 (Kyiv+Israel) 11:46 (UTC) 20:46 (Japan) 17:16 (India)"
$ws.Range("B94").Value = "***"
$ws.Range("C94").Value = "***"
$ws.Range("D94").Value = 3.881
$ws.Range("E94").Value = -2.743
$ws.Range("A95").Value = "11.02.2022 14:38 (Kyiv+Israel) 12:38 (UTC) 21:38 (Japan) 18:08 (India)"
$ws.Range("B95").Value = "***"
$ws.Range("C95").Value = "***"
$ws.Range("D95").Value = 4.187
$ws.Range("E95").Value = -3.049
$ws.Range("A96").Value = "11.07.2022 11:08 (Kyiv+Israel) 09:08 (UTC) 18:08 (Japan) 14:38 (India)"
$ws.Range("B96").Value = "***"
$ws.Range("C96").Value = "***"
$ws.Range("D96").Value = 2.116
$ws.Range("E96").Value = -0.9780000000000002
$ws.Range("A97").Value = "11.07.2022 11:10 (Kyiv+Israel) 09:10 (UTC) 18:10 (Japan) 14:40 (India)"
$ws.Range("B97").Value = 1.807
$ws.Range("C97").Value = -1.08
$ws.Range("D97").Value = "***"
$ws.Range("E97").Value = "***"

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A90").Value = "11.01.2022 17:33 (Kyiv+Israel) 15:33 (UTC) 00:33 (Japan) 21:03 (India)"
$ws.Range("B90").Value = "***"
$ws.Range("C90").Value = "***"
$ws.Range("D90").Value = 1.222
$ws.Range("E90").Value = -0.629
$ws.Range("A91").Value = "11.01.2022 18:04 (Kyiv+Israel) 16:04 (UTC) 01:04 (Japan) 21:34 (India)"
$ws.Range("B91").Value = "***"
$ws.Range("C91").Value = "***"
$ws.Range("D91").Value = 3.854
$ws.Range("E91").Value = -3.261
$ws.Range("A92").Value = "11.02.2022 13:49 (Kyiv+Israel) 11:49 (UTC) 20:49 (Japan) 17:19 (India)"
$ws.Range("B92").Value = "***"
$ws.Range("C92").Value = "***"
$ws.Range("D92").Value = 0.772
$ws.Range("E92").Value = -0.179
$ws.Range("A93").Value = "11.02.2022 14:41 (Kyiv+Israel) 12:41 (UTC) 21:41 (Japan) 18:11 (India)"
$ws.Range("B93").Value = "***"
$ws.Range("C93").Value = "***"
$ws.Range("D93").Value = 1.356
$ws.Range("E93").Value = -0.7630000000000001
$ws.Range("A94").Value = "11.07.2022 11:12 (Kyiv+Israel) 09:12 (UTC) 18:12 (Japan) 14:42 (India)"
$ws.Range("B94").Value = "***"
$ws.Range("C94").Value = "***"
$ws.Range("D94").Value = 0.958
$ws.Range("E94").Value = -0.365
$ws.Range("A95").Value = "11.07.2022 11:13 (Kyiv+Israel) 09:13 (UTC) 18:13 (Japan) 14:43 (India)"
$ws.Range("B95").Value = 0.825
$ws.Range("C95").Value = -0.245
$ws.Range("D95").Value = "***"
$ws.Range("E95").Value = "***"

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A83").Value = "11.01.2022 18:07 (Kyiv+Israel) 16:07 (UTC) 01:07 (Japan) 21:37 (India)"
$ws.Range("B83").Value = "***"
$ws.Range("C83").Value = "***"
$ws.Range("D83").Value = 3.305
$ws.Range("E83").Value = -1.721
$ws.Range("A84").Value = "11.02.2022 13:52 (Kyiv+Israel) 11:52 (UTC) 20:52 (Japan) 17:22 (India)"
$ws.Range("B84").Value = "***"
$ws.Range("C84").Value = "***"
$ws.Range("D84").Value = 4.095
$ws.Range("E84").Value = -2.511
$ws.Range("A85").Value = "11.02.2022 14:43 (Kyiv+Israel) 12:43 (UTC) 21:43 (Japan) 18:13 (India)"
$ws.Range("B85").Value = "***"
$ws.Range("C85").Value = "***"
$ws.Range("D85").Value = 6.016
$ws.Range("E85").Value = -4.432
$ws.Range("A86").Value = "11.07.2022 11:14 (Kyiv+Israel) 09:14 (UTC) 18:14 (Japan) 14:44 (India)"
$ws.Range("B86").Value = "***"
$ws.Range("C86").Value = "***"
$ws.Range("D86").Value = 1.585
$ws.Range("E86").Value = -0.0009999999999998899
$ws.Range("A87").Value = "11.07.2022 11:16 (Kyiv+Israel) 09:16 (UTC) 18:16 (Japan) 14:46 (India)"
$ws.Range("B87").Value = 2.066
$ws.Range("C87").Value = -0.9619999999999997
$ws.Range("D87").Value = "***"
$ws.Range("E87").Value = "***"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A96").Value = "11.02.2022 13:41 (Kyiv+Israel) 11:41 (UTC) 20:41 (Japan) 17:11 (India)"
$ws.Range("B96").Value = "***"
$ws.Range("C96").Value = "***"
$ws.Range("D96").Value = 3.199
$ws.Range("E96").Value = -2.277
$ws.Range("A97").Value = "11.02.2022 14:47 (Kyiv+Israel) 12:47 (UTC) 21:47 (Japan) 18:17 (India)"
$ws.Range("B97").Value = "***"
$ws.Range("C97").Value = "***"
$ws.Range("D97").Value = 1.434
$ws.Range("E97").Value = -0.5119999999999999
$ws.Range("A98").Value = "11.07.2022 11:17 (Kyiv+Israel) 09:17 (UTC) 18:17 (Japan) 14:47 (India)"
$ws.Range("B98").Value = "***"
$ws.Range("C98").Value = "***"
$ws.Range("D98").Value = 1.03
$ws.Range("E98").Value = -0.108
$ws.Range("A99").Value = "11.07.2022 11:24 (Kyiv+Israel) 09:24 (UTC) 18:24 (Japan) 14:54 (India)"
$ws.Range("B99").Value = 0.972
$ws.Range("C99").Value = -0.3069999999999999
$ws.Range("D99").Value = "***"
$ws.Range("E99").Value = "***"
